# Normalize the "Recorded By" column (G): entries recorded as
# "System, <email>" should instead read "<email>, System" -- i.e. put the
# actual user first and the automated "System" actor last. Entries that
# involve the backup/backdoor account are left exactly as-is, and any
# cell that isn't a simple two-part "System, <email>" value is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $text = $cell.Value2

    if ($text -eq $null) { continue }

    $parts = $text -split ','
    if ($parts.Count -eq 2) {
        $first = $parts[0].Trim()
        $second = $parts[1].Trim()

        if ($first -eq 'System' -and $second -notlike '*backdoor*') {
            $cell.Value = "$second, $first"
        }
    }
}
